$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("for Resp")

# Clear the "Remark" value in B1 (was "Testing for you you you")
$ws.Range("B1").Value = ""

# Clear the Tissue value in B8 (was "T-23000")
$ws.Range("B8").Value = ""

# B9 keeps "RJ" (string table renumbers after removed strings)
$ws.Range("B9").Value = "RJ"

# Remove stray numeric values
$ws.Range("B10").Value = ""
$ws.Range("B11").Value = ""

# Clear all the "Not Detected" result cells B12:B30
$ws.Range("B12:B30").Value = ""

# Update the selection to match the new active range
$ws.Activate()
$ws.Range("B12:B30").Select()
